# "edit index.html and fix style.css" -- the actual spreadsheet change:
# the "JENIS PENGADAAN BRG TERBANYAK" sheet collapses its Tahun/Jenis
# Pengadaan/Jumlah Pengadaan 3-column layout into a 2-column layout where
# the year + procurement type are merged into column A text, and column B
# holds the integer "Jumlah Pengadaan" count. That sheet also becomes the
# active tab (instead of "RATA RATA HARGA PENAWARAN").

$wb = $excel.ActiveWorkbook
$ws5 = $wb.Worksheets.Item("JENIS PENGADAAN BRG TERBANYAK")

# New width for column A (widened to fit the merged "YYYY (Jenis ...)" text).
# Set before clearing so the <cols> width entry survives the Cells.Clear() below.
$ws5.Range("A1").ColumnWidth = 34.8

# Wipe the old sheetData (A:C, 6 rows) clean -- this also drops the now
# unused column-C cell styling left behind by a plain ClearContents.
$ws5.Cells.Clear()

# Header row
$ws5.Range("A1").Value = "Tahun"
$ws5.Range("A1").Font.Bold = $true
$ws5.Range("B1").Value = "Jumlah Pengadaan"
$ws5.Range("B1").Font.Bold = $true

# Data rows: year+jenis merged into column A, count in column B (integer fmt)
$ws5.Range("A2").Value = "2017 (Pekerjaan Konstruksi)"
$ws5.Range("B2").Value = 512
$ws5.Range("B2").NumberFormat = "0"

$ws5.Range("A3").Value = "2018 (Pekerjaan Konstruksi"
$ws5.Range("B3").Value = 409
$ws5.Range("B3").NumberFormat = "0"

$ws5.Range("A4").Value = "2019 (Pengadaan Barang)"
$ws5.Range("B4").Value = 337
$ws5.Range("B4").NumberFormat = "0"

$ws5.Range("A5").Value = "2020 (Pengadaan Barang)"
$ws5.Range("B5").Value = 232
$ws5.Range("B5").NumberFormat = "0"

$ws5.Range("A6").Value = "2021 (Pekerjaan Konstruksi)"
$ws5.Range("B6").Value = 228
$ws5.Range("B6").NumberFormat = "0"

# Make this sheet the active tab/selection (was "RATA RATA HARGA PENAWARAN")
$ws5.Activate()
$ws5.Range("F11").Select()
